# Fruta / hortaliza, semanal
# Inserts a new weekly record at row 191 (pushing the previous rows 191-199
# down to 192-200) on the "Vega Modelo de Temuco" Pomelo price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 191, shifting existing rows down.
$ws.Rows.Item(191).Insert()

# Populate the new row with this week's data.
$ws.Cells.Item(191, 1).Value = 10
$ws.Cells.Item(191, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(191, 3).Value = "La Araucanía"
$ws.Cells.Item(191, 4).Value = 44615
$ws.Cells.Item(191, 5).Value = 9
$ws.Cells.Item(191, 6).Value = "Fruta"
$ws.Cells.Item(191, 7).Value = 100102
$ws.Cells.Item(191, 8).Value = "Cítricos"
$ws.Cells.Item(191, 9).Value = 100102006
$ws.Cells.Item(191, 10).Value = "Pomelo"
$ws.Cells.Item(191, 11).Value = "Start Ruby"
$ws.Cells.Item(191, 12).Value = "Primera"
$ws.Cells.Item(191, 13).Value = 80
$ws.Cells.Item(191, 14).Value = 15000
$ws.Cells.Item(191, 15).Value = 15000
$ws.Cells.Item(191, 16).Value = 15000
$ws.Cells.Item(191, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(191, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(191, 19).Value = 1000
$ws.Cells.Item(191, 20).Value = 15

# Match the style (date number format) used by the rest of column D.
$ws.Cells.Item(191, 4).NumberFormat = $ws.Cells.Item(192, 4).NumberFormat
